# Rename the client "Asanele Consultants (Pty) Ltd" -> "Comprac Energy (Pty) Ltd"
# everywhere it appears (document body + headers/footers of every section),
# and repoint the footer's picture "nice name" from shield.png to image.jpg
# (cosmetic pic:cNvPr/@name only - same embedded image/relationship).

$d = $word.ActiveDocument

$oldName = "Asanele Consultants (Pty) Ltd"
$newName = "Comprac Energy (Pty) Ltd"

# 1) Main document body / story.
$d.Content.Find.Execute($oldName, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newName, 2) | Out-Null

# 2) Every header/footer slot (primary=1, first page=2, even page=3) of
#    every section - Content.Find does not reach into these stories.
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)
    for ($hf = 1; $hf -le 3; $hf++) {
        $ftr = $sec.Footers.Item($hf)
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute($oldName, $true, $false, $false, $false, `
                                     $false, $true, 1, $false, $newName, 2) | Out-Null
        }
        $hdr = $sec.Headers.Item($hf)
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute($oldName, $true, $false, $false, $false, `
                                     $false, $true, 1, $false, $newName, 2) | Out-Null
        }
    }
}

# 3) The footer's picture internal name: pic:cNvPr/@name "shield.png" -> "image.jpg".
#    There is no InlineShape.Name property on the Word OM, so drop down to raw
#    OOXML: delete the shape, then InsertXML an identical drawing (same rId,
#    same extents/position, same paragraph formatting) with the new name, at
#    the exact spot the old shape occupied.
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)
    for ($hf = 1; $hf -le 3; $hf++) {
        $ftr = $sec.Footers.Item($hf)
        if (-not $ftr.Exists) { continue }
        $shapeCount = $ftr.Range.InlineShapes.Count
        for ($i = $shapeCount; $i -ge 1; $i--) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            if ($shp.Type -ne 3) { continue }   # wdInlineShapePicture

            $rng = $shp.Range
            $rng.Collapse(1)   # wdCollapseStart - keep a zero-length anchor

            $pPr = ""
            if ($shp.Range.Paragraphs.Item(1).Alignment -eq 1) {
                $pPr = "<w:pPr><w:jc w:val=""center""/></w:pPr>"
            }

            $xml = '<?xml version="1.0" standalone="yes"?>' + `
                   '<?mso-application progid="Word.Document"?>' + `
                   '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
                   '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
                   '<pkg:xmlData>' + `
                   '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
                   '<w:body><w:p>' + $pPr + '<w:r><w:drawing>' + `
                   '<wp:inline xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
                   '<wp:extent cx="360000" cy="447805"/>' + `
                   '<wp:docPr id="1" name="Picture 1"/>' + `
                   '<wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
                   '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
                   '<pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="image.jpg"/><pic:cNvPicPr/></pic:nvPicPr>' + `
                   '<pic:blipFill><a:blip r:embed="rId1"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
                   '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="360000" cy="447805"/></a:xfrm><a:prstGeom prst="rect"/></pic:spPr>' + `
                   '</pic:pic></a:graphicData></a:graphic></wp:inline>' + `
                   '</w:drawing></w:r></w:p></w:body></w:document>' + `
                   '</pkg:xmlData></pkg:part></pkg:package>'

            $shp.Delete()
            $rng.InsertXML($xml)
        }
    }
}
